$d = $word.ActiveDocument

$replacements = @(
    @("33×53=", "67×88="),
    @("22×63=", "41×84="),
    @("49×45=", "40×66="),
    @("89×14=", "59×60="),
    @("15×76=", "49×65="),
    @("33×58=", "29×98="),
    @("59×97=", "89×52="),
    @("74×55=", "99×57="),
    @("91×19=", "87×48="),
    @("46×70=", "29×74="),
    @("19×49=", "22×41="),
    @("67×87=", "96×26="),
    @("54×23=", "74×50="),
    @("58×75=", "33×37="),
    @("90×85=", "72×63="),
    @("99×49=", "71×41="),
    @("59×29=", "51×92="),
    @("43×46=", "38×18="),
    @("96×99=", "46×98="),
    @("11×12=", "87×17="),
    @("55×81=", "76×16="),
    @("99×86=", "49×75="),
    @("40×16=", "42×38="),
    @("29×68=", "85×53="),
    @("79×55=", "99×74=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
